$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (not numeric) interpretation for the Price column while we write,
# since several values look like plain numbers (e.g. "1.000", "240.66").
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '30.415.28'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').Value = '1.924.97'
$ws.Range('E3').Value = '  +4.02%  '
$ws.Range('D4').Value = '0.9996'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '240.66'
$ws.Range('E5').Value = '  +3.24%  '
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').Value = '0.4749'
$ws.Range('E7').Value = '  -0.33%  '
$ws.Range('D8').Value = '44.33'
$ws.Range('E8').Value = '  +2.29%  '
$ws.Range('D9').Value = '0.2856'
$ws.Range('E9').Value = '  +4.11%  '
$ws.Range('E10').Value = '  +3.82%  '
$ws.Range('D11').Value = '18.99'
$ws.Range('E11').Value = '  +7.96%  '
$ws.Range('D12').Value = '106.73'
$ws.Range('E12').Value = '  +26.08%  '
$ws.Range('D13').Value = '1.919.20'
$ws.Range('E13').Value = '  +3.74%  '
$ws.Range('D14').Value = '0.07613'
$ws.Range('E14').Value = '  +1.95%  '
$ws.Range('D15').Value = '5.124'
$ws.Range('E15').Value = '  +3.60%  '
$ws.Range('D16').Value = '0.6539'
$ws.Range('E16').Value = '  +4.83%  '
$ws.Range('D17').Value = '302.28'
$ws.Range('E17').Value = '  +23.03%  '
$ws.Range('D18').Value = '30.420.85'
$ws.Range('E18').Value = '  +0.11%  '
$ws.Range('D19').Value = '1.001'
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('E20').Value = '  +2.25%  '
$ws.Range('D21').Value = '2.181.74'
$ws.Range('E21').Value = '  +4.26%  '
$ws.Range('D22').Value = '0.000007481'
$ws.Range('E22').Value = '  +2.17%  '
$ws.Range('D23').Value = '5.291'
$ws.Range('E23').Value = '  +7.88%  '
$ws.Range('D24').Value = '1.0000'
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('D25').Value = '6.248'
$ws.Range('E25').Value = '  +5.91%  '
$ws.Range('D26').Value = '167.07'
$ws.Range('E26').Value = '  +1.48%  '
$ws.Range('D27').Value = '9.194'
$ws.Range('E27').Value = '  +1.29%  '
$ws.Range('D28').Value = '20.09'
$ws.Range('E28').Value = '  +11.67%  '
$ws.Range('E29').Value = '  +7.80%  '
$ws.Range('D30').Value = '0.1107'
$ws.Range('E30').Value = '  +7.71%  '
$ws.Range('D31').Value = '1.355'
$ws.Range('E31').Value = '  +0.25%  '
$ws.Range('D32').Value = '4.077'
$ws.Range('E32').Value = '  +1.02%  '
$ws.Range('D33').Value = '3.910'
$ws.Range('E33').Value = '  +2.51%  '
$ws.Range('D34').Value = '0.04978'
$ws.Range('E34').Value = '  +3.12%  '
$ws.Range('D35').Value = '0.7390'
$ws.Range('E35').Value = '  +6.40%  '
$ws.Range('E36').Value = '  +1.40%  '
$ws.Range('D37').Value = '2.749'
$ws.Range('E37').Value = '  +1.66%  '
$ws.Range('E38').Value = '  +1.98%  '
$ws.Range('D39').Value = '2.696'
$ws.Range('E39').Value = '  +0.50%  '
$ws.Range('E40').Value = '  +3.06%  '
$ws.Range('D41').Value = '0.8768'
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('D42').Value = '106.84'
$ws.Range('E42').Value = '  +0.11%  '
$ws.Range('D43').Value = '5.782'
$ws.Range('E43').Value = '  +5.10%  '
$ws.Range('D44').Value = '69.82'
$ws.Range('E44').Value = '  +10.93%  '
$ws.Range('D45').Value = '1.000'
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('D46').Value = '0.4126'
$ws.Range('E46').Value = '  +1.89%  '
$ws.Range('D47').Value = '7.239'
$ws.Range('E47').Value = '  +1.08%  '
$ws.Range('D48').Value = '9.247'
$ws.Range('E48').Value = '  +8.68%  '
$ws.Range('D49').Value = '34.80'
$ws.Range('E49').Value = '  +3.41%  '
$ws.Range('D50').Value = '0.1198'
$ws.Range('E50').Value = '  +0.27%  '
$ws.Range('D51').Value = '0.05618'
$ws.Range('E51').Value = '  +2.07%  '

# Restore the original (default) style now that the text values are committed.
$ws.Range("D2:D51").Style = "Normal"
